# "refining benford on pe"
# Update the computed cluster-summary statistics in columns B:K for rows 2-8
# (cluster ids 0-6) on the active sheet, reflecting a refined clustering run.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 2271.325892857143
$ws.Range("C2").Value = 224
$ws.Range("D2").Value = 0.6014461657843614
$ws.Range("E2").Value = 0.2539456642590378
$ws.Range("F2").Value = 0.06698943551544333
$ws.Range("G2").Value = 0.01174509978070864
$ws.Range("H2").Value = 0.1423715009126928
$ws.Range("I2").Value = 0.1140977635062004
$ws.Range("J2").Value = 0.008928571428571428

$ws.Range("B3").Value = 19539.40024630542
$ws.Range("C3").Value = 812
$ws.Range("D3").Value = 0.4350530495408141
$ws.Range("E3").Value = 0.2417765521100257
$ws.Range("F3").Value = 0.07536795491979049
$ws.Range("G3").Value = 0.01515534984687518
$ws.Range("H3").Value = 0.1408870074567197
$ws.Range("I3").Value = 0.1460742656532444
$ws.Range("J3").Value = 0.002463054187192118

$ws.Range("B4").Value = 831866.7717842323
$ws.Range("C4").Value = 245
$ws.Range("D4").Value = -0.1777495855196479
$ws.Range("E4").Value = 0.2475516863678265
$ws.Range("F4").Value = -0.008975337249607647
$ws.Range("G4").Value = 0.05573450551169317
$ws.Range("H4").Value = 0.2670605357714252
$ws.Range("I4").Value = 0.1851396145940806
$ws.Range("J4").Value = 0.07755102040816327

$ws.Range("B5").Value = 157672.180952381
$ws.Range("C5").Value = 424
$ws.Range("D5").Value = 0.1347852826689614
$ws.Range("E5").Value = 0.2489950878245923
$ws.Range("F5").Value = 0.01866252965682585
$ws.Range("G5").Value = 0.0452689156454677
$ws.Range("H5").Value = 0.2982780145768612
$ws.Range("I5").Value = 0.1665209188422025
$ws.Range("J5").Value = 0.05424528301886793

$ws.Range("B6").Value = 7683.529307282416
$ws.Range("C6").Value = 563
$ws.Range("D6").Value = 0.4325447714411564
$ws.Range("E6").Value = 0.296234995268464
$ws.Range("F6").Value = 0.07849168295753513
$ws.Range("G6").Value = 0.01247915707274186
$ws.Range("H6").Value = 0.1171674075059611
$ws.Range("I6").Value = 0.1422981312854583
$ws.Range("J6").Value = 0.01065719360568384
$ws.Range("K6").Value = 0.9982238010657194

$ws.Range("B7").Value = 35788.86301369863
$ws.Range("C7").Value = 219
$ws.Range("D7").Value = 0.366121662709347
$ws.Range("E7").Value = 0.2008177408075416
$ws.Range("F7").Value = 0.05780469874479671
$ws.Range("G7").Value = 0.02812309667813223
$ws.Range("H7").Value = 0.2294493837246292
$ws.Range("I7").Value = 0.1642803448217128
$ws.Range("J7").Value = 0.0182648401826484
$ws.Range("K7").Value = 1

$ws.Range("B8").Value = 51393.74363057325
$ws.Range("C8").Value = 628
$ws.Range("D8").Value = 0.257235473314821
$ws.Range("E8").Value = 0.2705787708868452
$ws.Range("F8").Value = 0.0507314244746194
$ws.Range("G8").Value = 0.03092498997316032
$ws.Range("H8").Value = 0.2154643038922539
$ws.Range("I8").Value = 0.1647149759918513
$ws.Range("J8").Value = 0.03503184713375796
